# "segunda actualizacion de datos" - append the latest risk-data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 used to be the last row (31-Dec data was not yet in); now fresh
# numbers come in for that date and it becomes a regular (non-last) row.
$ws.Range("A13").Value = 45657
$ws.Range("B13").Value = 0.030107279294914118
$ws.Range("C13").Value = 0.024543712565587263
$ws.Range("D13").Value = 0.012526613590292716

# Give row 13 the same "thick separator" row height used by every other
# interior data row (rows 1-12) so it visually matches them now that it's
# no longer the final row of the series.
$ws.Rows(13).RowHeight = 15

# New last row: 1-Jan data. Pull the cell formatting (date format / percent
# style) from an existing data row so the new cells match the rest of the
# column without minting new styles, then fill in the values.
$ws.Range("A2:D2").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A14").Value = 45658
$ws.Range("B14").Value = 0.027874542433205196
$ws.Range("C14").Value = 0.027110924433475184
$ws.Range("D14").Value = 0.011051310383238458

# Leave the current selection where it ended up after entering the new row.
[void]$ws.Range("G16").Select()
